$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

$ws.Range("H74").Value = 251941.5
$ws.Range("J74").Value = 500000
$ws.Range("L74").Value = 500000
$ws.Range("N74").Value = -501872

$ws.Range("H77").Value = 251941.5
$ws.Range("J77").Value = 500000
$ws.Range("L77").Value = 2500000
$ws.Range("N77").Value = -2509360

$ws.Range("H137").Value = 2070
$ws.Range("I137").Value = 975.3333
$ws.Range("K137").Value = 2925.9999
$ws.Range("M137").Value = -375.9998999999998

$ws.Range("H141").Value = 2394.1428
$ws.Range("I141").Value = 2394.1428
$ws.Range("K141").Value = 7182.428400000001
$ws.Range("M141").Value = -2002.428400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

$ws.Range("H132").Value = 1529.591
$ws.Range("I132").Value = 1425.1111
$ws.Range("K132").Value = 4275.3333
$ws.Range("M132").Value = -1745.3333

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 95
$ws.Range("I22").Value = 95
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 95
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 78
$ws.Range("N22").ClearContents()

$ws.Range("H80").Value = 700.8570999999999
$ws.Range("J80").Value = 703.4
$ws.Range("L80").Value = 703.4
$ws.Range("N80").Value = -2699.4

$ws.Range("H83").Value = 700.8570999999999
$ws.Range("J83").Value = 703.4
$ws.Range("L83").Value = 3517
$ws.Range("N83").Value = -13501

$ws.Range("H94").Value = 2269.75
$ws.Range("I94").Value = 2269.75
$ws.Range("K94").Value = 2269.75
$ws.Range("M94").Value = -1818.75

$ws.Range("H134").Value = 3688.7
$ws.Range("I134").Value = 3688.7
$ws.Range("K134").Value = 11066.1
$ws.Range("M134").Value = -8531.099999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1814.7018
$ws.Range("I31").Value = 1374.1666
$ws.Range("K31").Value = 1374.1666
$ws.Range("M31").Value = -1079.1666

$ws.Range("H34").Value = 1814.7018
$ws.Range("I34").Value = 1374.1666
$ws.Range("K34").Value = 1374.1666
$ws.Range("M34").Value = -1172.1666

$ws.Range("H58").Value = 5418.5264
$ws.Range("I58").Value = 4891.353
$ws.Range("K58").Value = 4891.353
$ws.Range("M58").Value = -4688.353

$ws.Range("H76").Value = 4999
$ws.Range("I76").Value = 4999
$ws.Range("K76").Value = 4999
$ws.Range("M76").Value = -4684

$ws.Range("H79").Value = 4999
$ws.Range("I79").Value = 4999
$ws.Range("K79").Value = 4999
$ws.Range("M79").Value = -3907

$ws.Range("H122").Value = 922.25
$ws.Range("I122").Value = 922.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2766.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -316.75
$ws.Range("N122").ClearContents()

$ws.Range("H136").Value = 5418.5264
$ws.Range("I136").Value = 4891.353
$ws.Range("K136").Value = 14674.059
$ws.Range("M136").Value = -12124.059

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 231.5
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 256.66666
$ws.Range("K2").Value = 30
$ws.Range("L2").Value = 1539.99996
$ws.Range("M2").Value = 83
$ws.Range("N2").Value = -1765.99996

$ws.Range("H5").Value = 1704.875
$ws.Range("I5").Value = 1748.7
$ws.Range("J5").Value = 1631.8334
$ws.Range("K5").Value = 5246.1
$ws.Range("L5").Value = 4895.5002
$ws.Range("M5").Value = -5134.1
$ws.Range("N5").Value = -5119.5002

$ws.Range("H7").Value = 362
$ws.Range("I7").Value = 263.25
$ws.Range("J7").Value = 427.83334
$ws.Range("K7").Value = 789.75
$ws.Range("L7").Value = 1283.50002
$ws.Range("M7").Value = -677.75
$ws.Range("N7").Value = -1507.50002

$ws.Range("H38").Value = 75.333336
$ws.Range("I38").Value = 45.5
$ws.Range("J38").Value = 135
$ws.Range("K38").Value = 136.5
$ws.Range("L38").Value = 405
$ws.Range("M38").Value = 210.5
$ws.Range("N38").Value = -1099

$ws.Range("H135").Value = 1704.875
$ws.Range("I135").Value = 1748.7
$ws.Range("J135").Value = 1631.8334
$ws.Range("K135").Value = 15738.3
$ws.Range("L135").Value = 14686.5006
$ws.Range("M135").Value = -13203.3
$ws.Range("N135").Value = -19756.5006

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 924.2
$ws.Range("I22").Value = 924.2
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 924.2
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -629.2
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 924.2
$ws.Range("I27").Value = 924.2
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 924.2
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -817.2
$ws.Range("N27").ClearContents()

$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

$ws.Range("H93").Value = 750
$ws.Range("I93").Value = 750
$ws.Range("K93").Value = 750
$ws.Range("M93").Value = 498

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws.Range("H136").Value = 3494
$ws.Range("I136").Value = 3494
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10482
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -7932
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 39000
$ws.Range("I38").Value = 39000
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 39000
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -38527
$ws.Range("N38").ClearContents()

$ws.Range("H136").Value = 6131.65
$ws.Range("I136").Value = 6305.3335
$ws.Range("J136").Value = 5610.6
$ws.Range("K136").Value = 18916.0005
$ws.Range("L136").Value = 16831.8
$ws.Range("M136").Value = -16366.0005
$ws.Range("N136").Value = -21931.8

